# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 36 (pushing the existing rows 36-43
# down to 37-44) and populate it with the latest "Zapallo italiano" price
# record for Comercializadora del Agro de Limarí (Coquimbo).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 36; existing rows 36:43 shift to 37:44.
$ws.Rows.Item(36).Insert()

$ws.Range("A36").Value = 2
$ws.Range("B36").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C36").Value = "Coquimbo"
$ws.Range("D36").Value = 44559
$ws.Range("E36").Value = 4
$ws.Range("F36").Value = 100112032
$ws.Range("G36").Value = "Zapallo italiano"
$ws.Range("H36").Value = "Sin especificar"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 560
$ws.Range("K36").Value = 6000
$ws.Range("L36").Value = 7000
$ws.Range("M36").Value = 6500
$ws.Range("N36").Value = "`$/caja 60 unidades"
$ws.Range("O36").Value = "Provincia de Limarí"
$ws.Range("P36").Value = 108
$ws.Range("Q36").Value = 60
$ws.Range("R36").Value = "Hortaliza"
